$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column mapping: D=4, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20

# Row 2
$ws.Cells.Item(2, 4).Value = 44357
$ws.Cells.Item(2, 13).Value = 10
$ws.Cells.Item(2, 14).Value = 38000
$ws.Cells.Item(2, 15).Value = 38000
$ws.Cells.Item(2, 16).Value = 38000
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(2, 18).Value = "Perú"
$ws.Cells.Item(2, 19).Value = 2111
$ws.Cells.Item(2, 20).Value = 18

# Row 3
$ws.Cells.Item(3, 4).Value = 44424
$ws.Cells.Item(3, 13).Value = 15
$ws.Cells.Item(3, 14).Value = 35000
$ws.Cells.Item(3, 15).Value = 35000
$ws.Cells.Item(3, 16).Value = 35000
$ws.Cells.Item(3, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(3, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 19).Value = 1944
$ws.Cells.Item(3, 20).Value = 18

# Row 4
$ws.Cells.Item(4, 4).Value = 44377
$ws.Cells.Item(4, 13).Value = 30
$ws.Cells.Item(4, 14).Value = 40000
$ws.Cells.Item(4, 15).Value = 40000
$ws.Cells.Item(4, 16).Value = 40000
$ws.Cells.Item(4, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(4, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 19).Value = 2222
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = 44364
$ws.Cells.Item(5, 13).Value = 90
$ws.Cells.Item(5, 14).Value = 1700
$ws.Cells.Item(5, 15).Value = 1700
$ws.Cells.Item(5, 16).Value = 1700
$ws.Cells.Item(5, 17).Value = "$/kilo"
$ws.Cells.Item(5, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 19).Value = 1700
$ws.Cells.Item(5, 20).Value = 1

# Row 6
$ws.Cells.Item(6, 4).Value = 44405
$ws.Cells.Item(6, 13).Value = 10
$ws.Cells.Item(6, 14).Value = 35000
$ws.Cells.Item(6, 15).Value = 35000
$ws.Cells.Item(6, 16).Value = 35000
$ws.Cells.Item(6, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(6, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 19).Value = 1944
$ws.Cells.Item(6, 20).Value = 18

# Row 7
$ws.Cells.Item(7, 4).Value = 44363
$ws.Cells.Item(7, 13).Value = 144
$ws.Cells.Item(7, 14).Value = 1700
$ws.Cells.Item(7, 15).Value = 1700
$ws.Cells.Item(7, 16).Value = 1700
$ws.Cells.Item(7, 17).Value = "$/kilo"
$ws.Cells.Item(7, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 19).Value = 1700
$ws.Cells.Item(7, 20).Value = 1

# Row 8
$ws.Cells.Item(8, 4).Value = 44431
$ws.Cells.Item(8, 13).Value = 30
$ws.Cells.Item(8, 14).Value = 35000
$ws.Cells.Item(8, 15).Value = 35000
$ws.Cells.Item(8, 16).Value = 35000
$ws.Cells.Item(8, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(8, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 19).Value = 1944
$ws.Cells.Item(8, 20).Value = 18

# Row 9
$ws.Cells.Item(9, 4).Value = 44418
$ws.Cells.Item(9, 13).Value = 30
$ws.Cells.Item(9, 14).Value = 35000
$ws.Cells.Item(9, 15).Value = 35000
$ws.Cells.Item(9, 16).Value = 35000
$ws.Cells.Item(9, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(9, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 19).Value = 1944
$ws.Cells.Item(9, 20).Value = 18

# Row 10
$ws.Cells.Item(10, 4).Value = 44392
$ws.Cells.Item(10, 13).Value = 20
$ws.Cells.Item(10, 14).Value = 35000
$ws.Cells.Item(10, 15).Value = 35000
$ws.Cells.Item(10, 16).Value = 35000
$ws.Cells.Item(10, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(10, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 19).Value = 1944
$ws.Cells.Item(10, 20).Value = 18

# Row 11
$ws.Cells.Item(11, 4).Value = 44433
$ws.Cells.Item(11, 13).Value = 15
$ws.Cells.Item(11, 14).Value = 35000
$ws.Cells.Item(11, 15).Value = 35000
$ws.Cells.Item(11, 16).Value = 35000
$ws.Cells.Item(11, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(11, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 19).Value = 1944
$ws.Cells.Item(11, 20).Value = 18

# Row 12
$ws.Cells.Item(12, 4).Value = 44434
$ws.Cells.Item(12, 13).Value = 40
$ws.Cells.Item(12, 14).Value = 35000
$ws.Cells.Item(12, 15).Value = 35000
$ws.Cells.Item(12, 16).Value = 35000
$ws.Cells.Item(12, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(12, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 19).Value = 1944
$ws.Cells.Item(12, 20).Value = 18

# Row 13
$ws.Cells.Item(13, 4).Value = 44264
$ws.Cells.Item(13, 13).Value = 20
$ws.Cells.Item(13, 14).Value = 40000
$ws.Cells.Item(13, 15).Value = 40000
$ws.Cells.Item(13, 16).Value = 40000
$ws.Cells.Item(13, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(13, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 19).Value = 2222
$ws.Cells.Item(13, 20).Value = 18

# Row 14
$ws.Cells.Item(14, 4).Value = 44432
$ws.Cells.Item(14, 13).Value = 10
$ws.Cells.Item(14, 14).Value = 35000
$ws.Cells.Item(14, 15).Value = 35000
$ws.Cells.Item(14, 16).Value = 35000
$ws.Cells.Item(14, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(14, 18).Value = "Perú"
$ws.Cells.Item(14, 19).Value = 1944
$ws.Cells.Item(14, 20).Value = 18

# Row 15
$ws.Cells.Item(15, 4).Value = 44369
$ws.Cells.Item(15, 13).Value = 5
$ws.Cells.Item(15, 14).Value = 35000
$ws.Cells.Item(15, 15).Value = 35000
$ws.Cells.Item(15, 16).Value = 35000
$ws.Cells.Item(15, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(15, 18).Value = "Perú"
$ws.Cells.Item(15, 19).Value = 1944
$ws.Cells.Item(15, 20).Value = 18

# Row 16
$ws.Cells.Item(16, 4).Value = 44294
$ws.Cells.Item(16, 13).Value = 15
$ws.Cells.Item(16, 14).Value = 35000
$ws.Cells.Item(16, 15).Value = 35000
$ws.Cells.Item(16, 16).Value = 35000
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(16, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 19).Value = 1944
$ws.Cells.Item(16, 20).Value = 18

# Row 17
$ws.Cells.Item(17, 4).Value = 44379
$ws.Cells.Item(17, 13).Value = 10
$ws.Cells.Item(17, 14).Value = 30000
$ws.Cells.Item(17, 15).Value = 30000
$ws.Cells.Item(17, 16).Value = 30000
$ws.Cells.Item(17, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(17, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 19).Value = 1667
$ws.Cells.Item(17, 20).Value = 18

# Row 18
$ws.Cells.Item(18, 4).Value = 44438
$ws.Cells.Item(18, 13).Value = 25
$ws.Cells.Item(18, 14).Value = 35000
$ws.Cells.Item(18, 15).Value = 35000
$ws.Cells.Item(18, 16).Value = 35000
$ws.Cells.Item(18, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(18, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(18, 19).Value = 1944
$ws.Cells.Item(18, 20).Value = 18

# Row 19
$ws.Cells.Item(19, 4).Value = 44442
$ws.Cells.Item(19, 13).Value = 15
$ws.Cells.Item(19, 14).Value = 35000
$ws.Cells.Item(19, 15).Value = 35000
$ws.Cells.Item(19, 16).Value = 35000
$ws.Cells.Item(19, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(19, 18).Value = "Perú"
$ws.Cells.Item(19, 19).Value = 1944
$ws.Cells.Item(19, 20).Value = 18

# Row 20
$ws.Cells.Item(20, 4).Value = 44435
$ws.Cells.Item(20, 13).Value = 10
$ws.Cells.Item(20, 14).Value = 35000
$ws.Cells.Item(20, 15).Value = 35000
$ws.Cells.Item(20, 16).Value = 35000
$ws.Cells.Item(20, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(20, 18).Value = "Perú"
$ws.Cells.Item(20, 19).Value = 1944
$ws.Cells.Item(20, 20).Value = 18

# Row 21
$ws.Cells.Item(21, 4).Value = 44435
$ws.Cells.Item(21, 13).Value = 105
$ws.Cells.Item(21, 14).Value = 35000
$ws.Cells.Item(21, 15).Value = 35000
$ws.Cells.Item(21, 16).Value = 35000
$ws.Cells.Item(21, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(21, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 19).Value = 1944
$ws.Cells.Item(21, 20).Value = 18

# Row 22
$ws.Cells.Item(22, 4).Value = 44279
$ws.Cells.Item(22, 13).Value = 30
$ws.Cells.Item(22, 14).Value = 35000
$ws.Cells.Item(22, 15).Value = 36000
$ws.Cells.Item(22, 16).Value = 35667
$ws.Cells.Item(22, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(22, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 19).Value = 1982
$ws.Cells.Item(22, 20).Value = 18

# New row 22: copy constant columns from row 21 (A,B,C,E,F,G,H,I,J,K,L are identical across all rows)
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(22, 3).Value = "La Araucanía"
$ws.Cells.Item(22, 5).Value = 9
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100108
$ws.Cells.Item(22, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(22, 9).Value = 100108003
$ws.Cells.Item(22, 10).Value = "Maracuyá"
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Primera"

# Apply date style (style index 2, numFmtId 165) to new D22 cell, matching D2:D21
$ws.Range("D22").NumberFormat = $ws.Range("D21").NumberFormat

